# "Cleanud up code, fixed some bugs"
#
# Row 25 ("Cleanup code") is being picked up: time spent is started at 0h,
# a note is added, and the status cell is recoloured from the old
# orange/"theme" fill to a new custom yellow-green fill (same colour the
# other in-progress rows like C22 use conceptually, but a freshly-picked
# custom colour here) to flag it as the task currently being worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Time spent so far on "Cleanup code" -> just started, 0 hours.
# Copy the time (h:mm) number format used by the other "Time Spent" cells.
$ws.Range("B25").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B25").Value2 = 0

# Give the status cell (currently "BUSY") a new custom highlight colour
# (FFB7C624) picked from the colour picker - replaces the old theme fill.
$ws.Range("C25").Interior.Color = 2410167

# Leave a note about the task.
$ws.Range("D25").Value2 = "I'll keep on cleaning until deadline"

# Cursor ended up on A26 after the edits.
$ws.Range("A26").Select()
